$wb = $excel.ActiveWorkbook

$wsOcc = $wb.Worksheets.Item("Occurrences")
$wsSampling = $wb.Worksheets.Item("Sampling Events")

# The "individualCount" (column M) values for the vegetation occurrence
# rows are not actually known from the source thesis, so replace the
# placeholder "1" with "?" for every data row (2-33).
$wsOcc.Range("M2:M33").Value = "?"

# Rows 34-36 were animal (Animalia / Spilocuscus maculatus / Phalanger
# orientalis / Phalanger permixtio) records mistakenly included in the
# vegetation occurrences sheet - remove them.
$wsOcc.Rows.Item(34).Resize(3).EntireRow.Delete()

# Reset the view: select H2 (also resets the scrolled-away topLeftCell),
# then make the "Sampling Events" sheet the active tab again.
$wsOcc.Activate()
$wsOcc.Range("H2").Select()

$wsSampling.Activate()

Write-Host "done"
